# feat: add 2022-Q3 data
#
# Inserts a new worksheet "2022-Q3" right after "总计" (Excel shuffles
# every later quarter sheet one slot to the right automatically; their
# names keep following their own data) and adds the corresponding
# summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet before the current second sheet
#    (today that is "2022-Q2").
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"
$q3.Cells.Clear()

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Style = "Normal"
}

$funds = @(
    @("515760", "华夏中证浙江国资创新发展ETF",   "2.04", "99.57", "3.65", "0.0745", 8),
    @("516530", "银华中证现代物流ETF",           "0.89", "97.53", "4.41", "0.0392", 6),
    @("516910", "富国中证现代物流ETF",           "0.78", "99.30", "4.48", "0.0349", 7),
    @("512190", "浙商汇金中证浙江凤凰行动50ETF", "0.48", "98.92", "3.41", "0.0164", 10),
    @("005120", "上投摩根量化多因子灵活配置混合", "0.19", "92.91", "1.78", "0.0034", 9),
    @("005966", "安信中证500指数增强C",          "0.16", "92.50", "1.07", "0.0017", 7),
    @("005965", "安信中证500指数增强A",          "0.10", "92.50", "1.07", "0.0011", 7)
)

for ($r = 0; $r -lt $funds.Length; $r++) {
    $fund = $funds[$r]
    $row = $r + 2

    # Column A: plain running index (0-based), numeric.
    $q3.Cells.Item($row, 1).Value = $r

    # Columns B-G are stored as text in this workbook (fund codes keep
    # their leading zeroes, and numeric-looking figures stay text
    # too) - force a text number format before writing the value,
    # then drop back to the unstyled "Normal" style so no stray
    # format survives on the cell.
    for ($c = 2; $c -le 7; $c++) {
        $target = $q3.Cells.Item($row, $c)
        $target.NumberFormat = "@"
        $target.Value = $fund[$c - 2]
        $target.Style = "Normal"
    }

    # Column H (rank) is a genuine number.
    $q3.Cells.Item($row, 8).Value = $fund[6]
}

# Re-apply the bold/boxed header style (column A index cells + header
# row) by copying it from a sheet that still carries it untouched.
$styleSource = $wb.Worksheets.Item("2022-Q2")
$styleSource.Range("A2").Copy()
for ($r = 2; $r -le 8; $r++) {
    $q3.Cells.Item($r, 1).PasteSpecial(-4122)
}
for ($r = 0; $r -lt $funds.Length; $r++) {
    $q3.Cells.Item($r + 2, 1).Value = $r
}

$styleSource.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. Add the new quarter's summary row to "总计" (row 2, pushing every
#    later quarter down by one and keeping the running index in
#    column A consistent).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

# Restore the index-column style on the new row (column A only).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.17

# Renumber the running index in column A (0-based row counter) for
# every quarter row that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
